$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date label
$ws.Range("A3").Value = "Date:30.05.19"

# Update sales figures for row 6 (Robiul)
$ws.Range("B6").Value = 158200
$ws.Range("C6").Value = 82
$ws.Range("D6").Value = 11

# Update sales figures for row 7 (Shohel)
$ws.Range("B7").Value = 108630
$ws.Range("C7").Value = 77
$ws.Range("D7").Value = 10

# Update sales figures for row 8 (Sodor)
$ws.Range("B8").Value = 115480
$ws.Range("C8").Value = 78

$excel.Calculate()
